# Actualización SmartScore desde Streamlit (Felis Bramley)
# Appends a new response row (row 14) to Sheet1 with the participant's
# submitted data. Columns holding numeric-looking SmartScore values are
# forced to Text via a leading apostrophe (matching the source export,
# which writes them as plain strings, not numbers), then the transient
# "quote prefix" style Excel applies is cleared back to Normal so the
# cells keep the sheet's default (unstyled) formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Felis Bramley_20251202_123358"

$ws.Range("B14").Value = "'"
$ws.Range("B14").Style = "Normal"

$ws.Range("C14").Value = "Felis Bramley"

$ws.Range("D14").Value = 18

$ws.Range("E14").Value = "Male"

$ws.Range("F14").Value = "2025-12-02 12:33:58"

$jsonG = @"
{
  "portion": 0.6,
  "diet": 0.5714285714285714,
  "salt": 0.2,
  "fat": 0.2,
  "natural": 0.4,
  "convenience": 0.4,
  "price": 0.0
}
"@
$ws.Range("G14").Value = $jsonG

$ws.Range("H14").Value = "Nongshim Neoguri Spicy Seafood"

$ws.Range("I14").Value = "'0.602"
$ws.Range("I14").Style = "Normal"

$ws.Range("J14").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Range("K14").Value = "Nissin Chow Mein Teriyaki Beef"

$ws.Range("L14").Value = "'0.494"
$ws.Range("L14").Style = "Normal"

$ws.Range("M14").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Range("N14").Value = "Nongshim Shin Ramyun"

$ws.Range("O14").Value = "'0.492"
$ws.Range("O14").Style = "Normal"

$ws.Range("P14").Value = "Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio"

$ws.Range("Q14").Value = "Amy’s Macaroni & Cheese (frozen)"

$ws.Range("R14").Value = "'0.776"
$ws.Range("R14").Style = "Normal"

$ws.Range("S14").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

$ws.Range("T14").Value = "Kraft Macaroni & Cheese Dinner"

$ws.Range("U14").Value = "'0.516"
$ws.Range("U14").Style = "Normal"

$ws.Range("V14").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Range("W14").Value = "Annie’s Shells & White Cheddar"

$ws.Range("X14").Value = "'0.511"
$ws.Range("X14").Style = "Normal"

$ws.Range("Y14").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

$ws.Range("Z14").Value = "Wild Planet Wild Tuna Pasta Salad"

$ws.Range("AA14").Value = "'0.808"
$ws.Range("AA14").Style = "Normal"

$ws.Range("AB14").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

$ws.Range("AC14").Value = "Kitchens of India Variety Pack"

$ws.Range("AD14").Value = "'0.524"
$ws.Range("AD14").Style = "Normal"

$ws.Range("AE14").Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"

$ws.Range("AF14").Value = "StarKist Chicken Creations (Chicken Salad)"

$ws.Range("AG14").Value = "'0.450"
$ws.Range("AG14").Style = "Normal"

$ws.Range("AH14").Value = "Portátil, saludable, fácil, buena textura, sabor suave"

# Re-fit the row height automatically (clears the explicit/custom height
# that got set implicitly when the multi-line G14 value was assigned),
# so the row keeps using the sheet's default automatic row height.
$ws.Rows.Item(14).EntireRow.AutoFit()
